$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Camote" (Vega Modelo de Temuco).
# It belongs right after the existing row 115 (chronologically among the
# other entries), so insert a new row at 116 and push everything else down.
$ws.Rows(116).Insert()

$newRow = 116
$ws.Cells.Item($newRow, 1).Value = 10
$ws.Cells.Item($newRow, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value = "2023-10-26"
$ws.Cells.Item($newRow, 5).Value = 9
$ws.Cells.Item($newRow, 6).Value = 100114002
$ws.Cells.Item($newRow, 7).Value = "Camote"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 80
$ws.Cells.Item($newRow, 11).Value = 24000
$ws.Cells.Item($newRow, 12).Value = 24000
$ws.Cells.Item($newRow, 13).Value = 24000
$ws.Cells.Item($newRow, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 1333
$ws.Cells.Item($newRow, 17).Value = 18
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
